$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A32").Value = "并音"
$ws.Range("B32").Value = 1

$ws.Range("A33").Value = "并音   多有力气"
$ws.Range("B33").Value = 1

$ws.Range("A34").Value = "语文幷"
$ws.Range("B34").Value = 1

$ws.Range("A35").Value = "拼音就能联续到我"
$ws.Range("B35").Value = 1

$ws.Range("A36").Value = "并音就能联续到我"
$ws.Range("B36").Value = 1

$wrapRange = $ws.Range("A35:A36")
$wrapRange.VerticalAlignment = -4108
$wrapRange.WrapText = $true

$ws.Range("A32").Select()
